# Atualizado por script em 01-11-2023 14:45
#
# This script re-applies the match data for the North Macedonia 1.MFL
# 2023-2024 sheet the way the scraping script would have: results for
# rows 8-9 (Aug 2023 round) and rows 60-71 (Oct 2023 rounds) are
# re-ordered, and a newly-played match (Struga 4-0 Vardar) is appended
# as row 72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Rows 8 & 9: the two matches swap position (only columns F:V -
# home/away teams, goals, odds, timestamps and url - move; the
# leading Indice/pais/torneio/temporada/data_partida columns A:E
# stay put).
# ---------------------------------------------------------------
$row8 = $ws.Range("F8:V8").Value()
$row9 = $ws.Range("F9:V9").Value()

$ws.Range("F8:V8").Value = $row9
$ws.Range("F9:V9").Value = $row8

# ---------------------------------------------------------------
# Rows 60-64: the five matches rotate down one slot, with the last
# one wrapping back around to the top (F:V only, A:E untouched).
# ---------------------------------------------------------------
$row60 = $ws.Range("F60:V60").Value()
$row61 = $ws.Range("F61:V61").Value()
$row62 = $ws.Range("F62:V62").Value()
$row63 = $ws.Range("F63:V63").Value()
$row64 = $ws.Range("F64:V64").Value()

$ws.Range("F60:V60").Value = $row64
$ws.Range("F61:V61").Value = $row60
$ws.Range("F62:V62").Value = $row61
$ws.Range("F63:V63").Value = $row62
$ws.Range("F64:V64").Value = $row63

# ---------------------------------------------------------------
# Rows 67-71: the five matches rotate by two slots (F:V only).
# ---------------------------------------------------------------
$row67 = $ws.Range("F67:V67").Value()
$row68 = $ws.Range("F68:V68").Value()
$row69 = $ws.Range("F69:V69").Value()
$row70 = $ws.Range("F70:V70").Value()
$row71 = $ws.Range("F71:V71").Value()

$ws.Range("F67:V67").Value = $row69
$ws.Range("F68:V68").Value = $row70
$ws.Range("F69:V69").Value = $row71
$ws.Range("F70:V70").Value = $row67
$ws.Range("F71:V71").Value = $row68

# ---------------------------------------------------------------
# New row 72: Struga 4-0 Vardar (31/10/2023), appended after row 71.
# Copy the formatting from row 71 first (style 1 on the Indice
# column, style 2 / date-time format on the data_partida column),
# then fill in the values.
# ---------------------------------------------------------------
$ws.Range("A71:V71").Copy()
$ws.Range("A72").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A72").Value = 71
$ws.Range("B72").Value = "north-macedonia"
$ws.Range("C72").Value = "1-mfl"
$ws.Range("D72").Value = "2023-2024"
$ws.Range("E72").Value = 45231.54166666666
$ws.Range("F72").Value = "Struga"
$ws.Range("G72").Value = 4
$ws.Range("H72").Value = "Vardar"
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 1.33
$ws.Range("K72").Value = "31/10/2023 01:12"
$ws.Range("L72").Value = 1.08
$ws.Range("M72").Value = "01/11/2023 12:58"
$ws.Range("N72").Value = 3.91
$ws.Range("O72").Value = "31/10/2023 01:12"
$ws.Range("P72").Value = 7.98
$ws.Range("Q72").Value = "01/11/2023 12:59"
$ws.Range("R72").Value = 6.9
$ws.Range("S72").Value = "31/10/2023 01:12"
$ws.Range("T72").Value = 27.05
$ws.Range("U72").Value = "01/11/2023 12:59"
$ws.Range("V72").Value = "https://www.betexplorer.com/football/north-macedonia/1-mfl/struga-vardar/QwlA9Dtc/"
